$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 311.57144
$ws.Range("I33").Value = 274
$ws.Range("K33").Value = 274
$ws.Range("M33").Value = -45
$ws.Range("H100").Value = 2653.4614
$ws.Range("I100").Value = 2324.375
$ws.Range("J100").Value = 3180
$ws.Range("K100").Value = 2324.375
$ws.Range("L100").Value = 3180
$ws.Range("M100").Value = -1783.375
$ws.Range("N100").Value = -4262

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2195.3
$ws.Range("I45").Value = 1513.4333
$ws.Range("J45").Value = 3218.1
$ws.Range("K45").Value = 1513.4333
$ws.Range("L45").Value = 3218.1
$ws.Range("M45").Value = -1136.4333
$ws.Range("N45").Value = -3972.1
$ws.Range("H74").Value = 2760.5833
$ws.Range("I74").Value = 2190.3333
$ws.Range("K74").Value = 2190.3333
$ws.Range("M74").Value = -1316.3333
$ws.Range("H77").Value = 2760.5833
$ws.Range("I77").Value = 2190.3333
$ws.Range("K77").Value = 10951.6665
$ws.Range("M77").Value = -6583.666499999999
$ws.Range("H97").Value = 3265
$ws.Range("I97").Value = 3152.5
$ws.Range("J97").Value = 3490
$ws.Range("K97").Value = 3152.5
$ws.Range("L97").Value = 3490
$ws.Range("M97").Value = -2656.5
$ws.Range("N97").Value = -4482
$ws.Range("H102").Value = 1221.1111
$ws.Range("I102").Value = 1305
$ws.Range("J102").Value = 550
$ws.Range("K102").Value = 1305
$ws.Range("L102").Value = 550
$ws.Range("M102").Value = 317
$ws.Range("N102").Value = -3794
$ws.Range("H110").Value = 4395.857
$ws.Range("I110").Value = 4447.75
$ws.Range("K110").Value = 4447.75
$ws.Range("M110").Value = -2402.75
$ws.Range("H115").Value = 39684
$ws.Range("J115").Value = 39684
$ws.Range("L115").Value = 39684
$ws.Range("N115").Value = -42818
$ws.Range("H139").Value = 50667.5
$ws.Range("J139").Value = 50667.5
$ws.Range("L139").Value = 50667.5
$ws.Range("N139").Value = -60947.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 409.7143
$ws.Range("I22").Value = 444.66666
$ws.Range("J22").Value = 200
$ws.Range("K22").Value = 444.66666
$ws.Range("L22").Value = 200
$ws.Range("M22").Value = -271.66666
$ws.Range("N22").Value = -546
$ws.Range("H94").Value = 2358.7222
$ws.Range("I94").Value = 1703.9166
$ws.Range("J94").Value = 3668.3333
$ws.Range("K94").Value = 1703.9166
$ws.Range("L94").Value = 3668.3333
$ws.Range("M94").Value = -1252.9166
$ws.Range("N94").Value = -4570.3333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 16366.482
$ws.Range("J31").Value = 5221.933
$ws.Range("L31").Value = 5221.933
$ws.Range("N31").Value = -5811.933
$ws.Range("H33").Value = 10779.8
$ws.Range("I33").Value = 4724.75
$ws.Range("J33").Value = 35000
$ws.Range("K33").Value = 4724.75
$ws.Range("L33").Value = 35000
$ws.Range("M33").Value = -4345.75
$ws.Range("N33").Value = -35758
$ws.Range("H34").Value = 16366.482
$ws.Range("J34").Value = 5221.933
$ws.Range("L34").Value = 5221.933
$ws.Range("N34").Value = -5625.933
$ws.Range("H68").Value = 69990
$ws.Range("J68").Value = 69990
$ws.Range("L68").Value = 69990
$ws.Range("N68").Value = -71488
$ws.Range("H71").Value = 69990
$ws.Range("J71").Value = 69990
$ws.Range("L71").Value = 209970
$ws.Range("N71").Value = -217458
$ws.Range("H105").Value = 25000504
$ws.Range("I105").Value = 31250376
$ws.Range("J105").Value = 1011
$ws.Range("K105").Value = 31250376
$ws.Range("L105").Value = 1011
$ws.Range("M105").Value = -31248629
$ws.Range("N105").Value = -4505
$ws.Range("H134").Value = 936.4
$ws.Range("J134").Value = 1418.5
$ws.Range("L134").Value = 4255.5
$ws.Range("N134").Value = -9325.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 751.15
$ws.Range("J131").Value = 751.15
$ws.Range("L131").Value = 2253.45
$ws.Range("N131").Value = -12333.45

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 26933.334
$ws.Range("J39").Value = 26933.334
$ws.Range("L39").Value = 26933.334
$ws.Range("N39").Value = -27997.334
$ws.Range("H126").Value = 4036.375
$ws.Range("I126").Value = 3383.3333
$ws.Range("J126").Value = 4612.5884
$ws.Range("K126").Value = 10149.9999
$ws.Range("L126").Value = 13837.7652
$ws.Range("M126").Value = -7679.999899999999
$ws.Range("N126").Value = -18777.7652

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 1000035.3
$ws.Range("I2").Value = 1100020
$ws.Range("J2").Value = 250150
$ws.Range("K2").Value = 1100020
$ws.Range("L2").Value = 250150
$ws.Range("M2").Value = -1099908
$ws.Range("N2").Value = -250374
$ws.Range("H36").Value = 34500
$ws.Range("J36").Value = 34500
$ws.Range("L36").Value = 34500
$ws.Range("N36").Value = -35624
$ws.Range("H68").Value = 3099.923
$ws.Range("I68").Value = 3144.3333
$ws.Range("J68").Value = 3000
$ws.Range("K68").Value = 3144.3333
$ws.Range("L68").Value = 3000
$ws.Range("M68").Value = -2395.3333
$ws.Range("N68").Value = -4498
$ws.Range("H71").Value = 3099.923
$ws.Range("I71").Value = 3144.3333
$ws.Range("J71").Value = 3000
$ws.Range("K71").Value = 15721.6665
$ws.Range("L71").Value = 15000
$ws.Range("M71").Value = -11977.6665
$ws.Range("N71").Value = -22488
$ws.Range("H82").Value = 2854.5454
$ws.Range("I82").Value = 4233.3335
$ws.Range("J82").Value = 1200
$ws.Range("K82").Value = 4233.3335
$ws.Range("L82").Value = 1200
$ws.Range("M82").Value = -3872.3335
$ws.Range("N82").Value = -1922
$ws.Range("H85").Value = 2854.5454
$ws.Range("I85").Value = 4233.3335
$ws.Range("J85").Value = 1200
$ws.Range("K85").Value = 4233.3335
$ws.Range("L85").Value = 1200
$ws.Range("M85").Value = -2985.3335
$ws.Range("N85").Value = -3696
$ws.Range("H136").Value = 36283.535
$ws.Range("I136").Value = 47386.637
$ws.Range("K136").Value = 142159.911
$ws.Range("M136").Value = -139609.911

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 949909.75
$ws.Range("I136").Value = 1403338.6
$ws.Range("K136").Value = 4210015.800000001
$ws.Range("M136").Value = -4207465.800000001

